$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new weekly price record was inserted for "Berenjena" (Región de Arica y
# Parinacota) at row 245, pushing the existing rows 245:294 down to 246:295.
$ws.Rows.Item(245).Insert()

$ws.Cells.Item(245, 1).Value = 9
$ws.Cells.Item(245, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(245, 3).Value = "Metropolitana"
$ws.Cells.Item(245, 4).Value = 44782
$ws.Cells.Item(245, 5).Value = 13
$ws.Cells.Item(245, 6).Value = 100112001
$ws.Cells.Item(245, 7).Value = "Berenjena"
$ws.Cells.Item(245, 8).Value = "Sin especificar"
$ws.Cells.Item(245, 9).Value = "Primera"
$ws.Cells.Item(245, 10).Value = 43
$ws.Cells.Item(245, 11).Value = 9000
$ws.Cells.Item(245, 12).Value = 10000
$ws.Cells.Item(245, 13).Value = 9512
$ws.Cells.Item(245, 14).Value = "`$/caja 50 unidades"
$ws.Cells.Item(245, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(245, 16).Value = 190
$ws.Cells.Item(245, 17).Value = 50
$ws.Cells.Item(245, 18).Value = "Hortaliza"
